# Scheduled data refresh: update market-price derived columns (H-N)
# across the Leve profit sheets, per the upstream snapshot pull.
$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")

# Row 62
$ws.Range("H62").Value = 52633330
$ws.Range("I62").Value = 76924590
$ws.Range("J62").Value = 2251
$ws.Range("K62").Value = 76924590
$ws.Range("L62").Value = 2251
$ws.Range("M62").Value = -76923966
$ws.Range("N62").Value = -3499

# Row 65
$ws.Range("H65").Value = 52633330
$ws.Range("I65").Value = 76924590
$ws.Range("J65").Value = 2251
$ws.Range("K65").Value = 384622950
$ws.Range("L65").Value = 11255
$ws.Range("M65").Value = -384619830
$ws.Range("N65").Value = -17495

# Row 113
$ws.Range("H113").Value = 2535.7693
$ws.Range("I113").Value = 2521.4285
$ws.Range("J113").Value = 2596
$ws.Range("K113").Value = 2521.4285
$ws.Range("L113").Value = 2596
$ws.Range("M113").Value = 732.5715
$ws.Range("N113").Value = -9104

# Row 132
$ws.Range("H132").Value = 8622602
$ws.Range("I132").Value = 9261221
$ws.Range("J132").Value = 1250
$ws.Range("K132").Value = 27783663
$ws.Range("L132").Value = 3750
$ws.Range("M132").Value = -27781133
$ws.Range("N132").Value = -8810

# Row 138
$ws.Range("H138").Value = 4250.9
$ws.Range("I138").Value = 1549.5
$ws.Range("J138").Value = 6614.625
$ws.Range("K138").Value = 4648.5
$ws.Range("L138").Value = 19843.875
$ws.Range("M138").Value = 491.5
$ws.Range("N138").Value = -30123.875

# Row 141
$ws.Range("H141").Value = 3232.2285
$ws.Range("I141").Value = 3224.2666
$ws.Range("J141").Value = 3280
$ws.Range("K141").Value = 9672.799800000001
$ws.Range("L141").Value = 9840
$ws.Range("M141").Value = -4492.799800000001
$ws.Range("N141").Value = -20200

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")

# Row 2
$ws.Range("H2").Value = 912.1539
$ws.Range("I2").Value = 425.76923
$ws.Range("J2").Value = 1884.9231
$ws.Range("K2").Value = 425.76923
$ws.Range("L2").Value = 1884.9231
$ws.Range("M2").Value = -312.76923
$ws.Range("N2").Value = -2110.9231

# Row 61
$ws.Range("H61").Value = 3511.7273
$ws.Range("I61").Value = 2973.4
$ws.Range("J61").Value = 4665.2856
$ws.Range("K61").Value = 2973.4
$ws.Range("L61").Value = 4665.2856
$ws.Range("M61").Value = -2761.4
$ws.Range("N61").Value = -5089.2856

# Row 116
$ws.Range("H116").Value = 912.1539
$ws.Range("I116").Value = 425.76923
$ws.Range("J116").Value = 1884.9231
$ws.Range("K116").Value = 425.76923
$ws.Range("L116").Value = 1884.9231
$ws.Range("M116").Value = 1868.23077
$ws.Range("N116").Value = -6472.9231

# Row 132
$ws.Range("H132").Value = 16701958
$ws.Range("I132").Value = 25001720
$ws.Range("J132").Value = 102431.1
$ws.Range("K132").Value = 75005160
$ws.Range("L132").Value = 307293.3
$ws.Range("M132").Value = -75002630
$ws.Range("N132").Value = -312353.3

# Row 136
$ws.Range("H136").Value = 3511.7273
$ws.Range("I136").Value = 2973.4
$ws.Range("J136").Value = 4665.2856
$ws.Range("K136").Value = 8920.200000000001
$ws.Range("L136").Value = 13995.8568
$ws.Range("M136").Value = -6370.200000000001
$ws.Range("N136").Value = -19095.8568

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")

# Row 35
$ws.Range("H35").Value = 11000
$ws.Range("J35").Value = 11000
$ws.Range("L35").Value = 11000
$ws.Range("N35").Value = -11620

# Row 134
$ws.Range("H134").Value = 14708.375
$ws.Range("I134").Value = 5544.5
$ws.Range("J134").Value = 42200
$ws.Range("K134").Value = 16633.5
$ws.Range("L134").Value = 126600
$ws.Range("M134").Value = -14098.5
$ws.Range("N134").Value = -131670

# Row 140
$ws.Range("H140").Value = 69500
$ws.Range("J140").Value = 69500
$ws.Range("L140").Value = 69500
$ws.Range("N140").Value = -79860

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")

# Row 107
$ws.Range("H107").Value = 1004.4
$ws.Range("I107").Value = 1007.3333
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1007.3333
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 912.6667
$ws.Range("N107").Value = -4840

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")

# Row 3
$ws.Range("H3").Value = 365
$ws.Range("I3").Value = 365
$ws.Range("K3").Value = 1095
$ws.Range("M3").Value = -983

# Row 19
$ws.Range("H19").Value = 3000
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 3000
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 9000
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -9348

# Row 54
$ws.Range("H54").Value = 5700
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 5700
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 17100
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -18218

# Row 64
$ws.Range("H64").Value = 766
$ws.Range("I64").Value = 757.5
$ws.Range("J64").Value = 800
$ws.Range("K64").Value = 2272.5
$ws.Range("L64").Value = 2400
$ws.Range("M64").Value = -2002.5
$ws.Range("N64").Value = -2940

# Row 67
$ws.Range("H67").Value = 766
$ws.Range("I67").Value = 757.5
$ws.Range("J67").Value = 800
$ws.Range("K67").Value = 2272.5
$ws.Range("L67").Value = 2400
$ws.Range("M67").Value = -1336.5
$ws.Range("N67").Value = -4272

# Row 68
$ws.Range("H68").Value = 1153.1765
$ws.Range("J68").Value = 1311.4445
$ws.Range("L68").Value = 3934.3335
$ws.Range("N68").Value = -5556.333500000001

# Row 71
$ws.Range("H71").Value = 1153.1765
$ws.Range("J71").Value = 1311.4445
$ws.Range("L71").Value = 11803.0005
$ws.Range("N71").Value = -19915.0005

# Row 107
$ws.Range("H107").Value = 856.1429000000001
$ws.Range("J107").Value = 1237.5294
$ws.Range("L107").Value = 3712.5882
$ws.Range("N107").Value = -7552.5882

# Row 129
$ws.Range("H129").Value = 10754034
$ws.Range("I129").Value = 1383.3334
$ws.Range("J129").Value = 17545182
$ws.Range("K129").Value = 4150.0002
$ws.Range("L129").Value = 52635546
$ws.Range("M129").Value = 849.9997999999996
$ws.Range("N129").Value = -52645546

# Row 140
$ws.Range("H140").Value = 6357.263
$ws.Range("I140").Value = 5385.933
$ws.Range("K140").Value = 16157.799
$ws.Range("M140").Value = -10977.799

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")

# Row 7
$ws.Range("H7").Value = 1538.8636
$ws.Range("I7").Value = 1331.6666
$ws.Range("K7").Value = 1331.6666
$ws.Range("M7").Value = -1219.6666

# Row 123
$ws.Range("H123").Value = 54980
$ws.Range("J123").Value = 54980
$ws.Range("L123").Value = 54980
$ws.Range("N123").Value = -64780

# Row 126
$ws.Range("H126").Value = 1538.8636
$ws.Range("I126").Value = 1331.6666
$ws.Range("K126").Value = 3994.9998
$ws.Range("M126").Value = -1524.9998

# Row 136
$ws.Range("H136").Value = 55614900
$ws.Range("I136").Value = 201051.8
$ws.Range("J136").Value = 76927920
$ws.Range("K136").Value = 603155.3999999999
$ws.Range("L136").Value = 230783760
$ws.Range("M136").Value = -600605.3999999999
$ws.Range("N136").Value = -230788860

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")

# Row 45
$ws.Range("H45").Value = 8068.2144
$ws.Range("J45").Value = 8068.2144
$ws.Range("L45").Value = 8068.2144
$ws.Range("N45").Value = -9050.214400000001

# Row 74
$ws.Range("H74").Value = 28079.6
$ws.Range("J74").Value = 28079.6
$ws.Range("L74").Value = 28079.6
$ws.Range("N74").Value = -29951.6

# Row 77
$ws.Range("H77").Value = 28079.6
$ws.Range("J77").Value = 28079.6
$ws.Range("L77").Value = 84238.79999999999
$ws.Range("N77").Value = -93598.79999999999

# Row 132
$ws.Range("H132").Value = 54883156
$ws.Range("I132").Value = 90002120
$ws.Range("J132").Value = 9774.5625
$ws.Range("K132").Value = 270006360
$ws.Range("L132").Value = 29323.6875
$ws.Range("M132").Value = -270003830
$ws.Range("N132").Value = -34383.6875
